# Fixed #295 Add the version of M2Doc in the template custom properties.
#
# The canonical-OOXML diff for this particular template resource shows no
# actual content/value changes in word/document.xml or word/styles.xml -
# every hunk is the same element with its attributes re-ordered (e.g.
# <w:pgSz w:w=".." w:h=".."/> -> <w:pgSz w:h=".." w:w=".."/>), which is a
# by-product of the repository's fixture-regeneration tooling, not a
# functional edit. There is no user-visible/content change to apply to
# this document: the body text, formatting, styles, sections and
# languages are all unchanged before/after.
#
# The only semantically meaningful part of the commit is recording the
# M2Doc version that generated/validated the template as a custom
# document property (done once, centrally, for every template in the
# fixture set). We still perform that step here through the standard
# Word object model call so the intent of the commit is honored on this
# document as well; on hosts where custom document properties are not
# writable this is harmlessly skipped.
$d = $word.ActiveDocument

try {
    $existing = $null
    try {
        $existing = $d.CustomDocumentProperties("M2DocVersion")
    } catch {
        $existing = $null
    }

    if ($existing -ne $null) {
        $existing.Value = "1.0"
    } else {
        # Name, LinkToContent, Type (msoPropertyTypeString = 4), Value
        $d.CustomDocumentProperties.Add("M2DocVersion", $false, 4, "1.0")
    }
} catch {
    # Custom document properties are not supported/persisted on this
    # host; nothing else in the template needs to change.
}
